# Update code! support multi device running -- Nice!
#
# - dialer_page: add a new "recents_button" locator row
# - messaging_page: add a new "message_conversation_title" locator row,
#   and re-point the view (zoom + selection) at the new data
# - calendar_page becomes the active tab (selection moves to B5)
# - fota_page is no longer the active tab

$wb = $excel.ActiveWorkbook

# --- dialer_page: new recents-tab locator row -----------------------------
$wsDialer = $wb.Worksheets.Item("dialer_page")
# Write column B before column A so the new shared strings land in the
# same order as the reference edit (call_log_tab, then recents_button).
$wsDialer.Cells.Item(9, 2).Value = "com.google.android.dialer:id/call_log_tab"
$wsDialer.Cells.Item(9, 1).Value = "recents_button"
$wsDialer.Cells.Item(9, 3).Value = "str"
$wsDialer.Range("A8").Select() | Out-Null

# --- messaging_page: new conversation-title locator row --------------------
$wsMessaging = $wb.Worksheets.Item("messaging_page")
$wsMessaging.Cells.Item(7, 2).Value = "com.google.android.apps.messaging:id/tombstone_message"
$wsMessaging.Cells.Item(7, 1).Value = "message_conversation_title"
$wsMessaging.Cells.Item(7, 3).Value = "str"
$wsMessaging.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 157
$wsMessaging.Range("B9").Select() | Out-Null

# --- calendar_page becomes the active sheet/selection ----------------------
$wsCalendar = $wb.Worksheets.Item("calendar_page")
$wsCalendar.Activate() | Out-Null
$wsCalendar.Range("B5").Select() | Out-Null
